$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the three runs that make up the middle of the "Postavljanje
#    kostura ..." bullet back into one run. The text does not change, only
#    the run split does, so re-"finding & replacing" the exact text that
#    spans those three runs with itself causes Word to normalise them into
#    a single run (matching the target diff) while leaving the surrounding
#    runs (the opening "(Antonela Bogdanic / 2 sata) ..." run and the
#    trailing run that carries the lastRenderedPageBreak) untouched.
# ---------------------------------------------------------------------------
$mid = "Primjećeni problemi i moguća poboljšanja"
$d.Content.Find.Execute($mid, $true, $false, $false, $false, $false, $true, 1, $false, $mid, 2)

# ---------------------------------------------------------------------------
# 2) Add a new "Dnevnik rada" bullet describing the player-name / move
#    counter work, right after the "Postavljanje kostura ..." bullet.
# ---------------------------------------------------------------------------
$kosturText = "(Antonela Bogdanić / 2 sata) Postavljanje kostura projekta i postavljanje na git, Uvod i Primjećeni problemi i moguća poboljšanja u ovom dokumentu, kostur ostalog dijela, početne preinake u programu, povećanje dimenzija, dodane nove varijable koje će nam trebati, izmjena završnog zaslona, uređivanje fotografija koje će biti završna pozadina ovisno o ishodu igre. "

$kosturIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($kosturText + "`r")) {
        $kosturIdx = $i
    }
}

$kosturPara = $d.Paragraphs.Item($kosturIdx)
$kosturPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($kosturIdx + 1)
$newPara1.Range.Text = "(Antonela Bogdanić / 1 sat) Dodavanje imena igrača i broja poteza na igru, uređivanje oblika i postavljanje da se mijenja boja igrača u crveno u ovisnosti o tome tko je na potezu te brojanje poteza. Dodana pogreška o potezu, koja ostoje sve dok se ne odigra ispravan potez."

# ---------------------------------------------------------------------------
# 3) Add two new bullets at the end of the "MOJE NAPOMENE:" list: one about
#    limiting player-name length, and a trailing empty bullet.
# ---------------------------------------------------------------------------
$lastText = "Mozda dodati u dokumentaciju ovu da smo htjeli snimati partije"

$lastIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($lastText + "`r")) {
        $lastIdx = $i
    }
}

$lastPara = $d.Paragraphs.Item($lastIdx)
$lastPara.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($lastIdx + 1)
$newPara2.Range.Text = "Dodat da imena budu manja od 10 znakova da se nebi desilo krivo kod teksta"

$newPara2.Range.InsertParagraphAfter()
